# [ANV] updating decay chain spreadsheet
#
# Adds a new "Cu Target Fractions " worksheet (copper target, analogous to
# the existing "HDPE Target Fractions" sheet), placed after "HDPE Target
# Fractions", and updates a couple of view-state bits (selection) on other
# sheets.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the selection/cursor on "HDPE Density" before moving on.
# ------------------------------------------------------------------
$wsHdpeDensity = $wb.Worksheets.Item("HDPE Density")
$wsHdpeDensity.Activate()
$wsHdpeDensity.Range("E24").Select()

# ------------------------------------------------------------------
# 2. Create the new "Cu Target Fractions " sheet by copying the
#    existing "HDPE Target Fractions" sheet (keeps identical column
#    widths / cell styles / number formats), placed right after it.
# ------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("HDPE Target Fractions")
$srcSheet.Copy([System.Reflection.Missing]::Value, $srcSheet)
$ws = $wb.Worksheets.Item($srcSheet.Index + 1)
$ws.Name = "Cu Target Fractions "

# --- Row 1 (headers): drop the extra "Middle (frac): Hang's
#     Normalization" header in column I -- the Cu sheet has no column I.
$ws.Range("I1").Clear()

# --- Row 2: Z = 29 (copper), fraction of atoms = 1.
$ws.Range("A2").Value = 29
$ws.Range("B2").Value = 1
$ws.Range("A3").Clear()

# add the (empty, but formatted) trailing cells on row 2
$ws.Range("F2").NumberFormat = "0.000E+00"
$ws.Range("G2").NumberFormat = "0.000E+00"
$ws.Range("H2").NumberFormat = "0.000E+00"
$ws.Range("I2").NumberFormat = "0.000E+00"

# --- Row 3: Cu-63 isotope.
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 63
$ws.Range("D3").Formula = "=69.15*(1-0.15)"
$ws.Range("E3").Formula = "=69.15*(1+0.15)"
$ws.Range("F3").Formula = "=B3*(D3+E3)/200"
$ws.Range("G3").Formula = "=B3*D3/100"
$ws.Range("H3").Formula = "=B3*E3/100"
$ws.Range("I3").NumberFormat = "0.000E+00"

# --- Row 4: Cu-65 isotope.
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 65
$ws.Range("D4").Formula = "=30.85*(0.85)"
$ws.Range("E4").Formula = "=30.85*1.15"
# F4/G4/H4 already carry the right relative formulas from the copied sheet
# (=B4*(D4+E4)/200, =B4*D4/100, =B4*E4/100) -- leave them as-is.
$ws.Range("I4").ClearContents()

# --- Row 5: total fraction of atoms is now just B2 (single element).
$ws.Range("B5").Formula = "=B2"

# --- Final selection / active-sheet state for the new sheet.
$ws.Range("F15").Select()
